$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("D3", "D4", "D6", "D7", "C10", "E14", "E15", "C16", "C18", "D18", "C19", "D19", "C20", "D20", "C21", "D21")

foreach ($cell in $cells) {
    $ws.Range($cell).Value = "-"
}
